$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Group = 15)
$ws.Range("B2").Value = -0.999490387077816
$ws.Range("C2").Value = -167965926.03901
$ws.Range("D2").Value = 0.956621583642287
$ws.Range("E2").Value = 0.999489950006617
$ws.Range("F2").Value = -0.956629205321652
$ws.Range("G2").Value = 123398.51151122
$ws.Range("H2").Value = 7053891.21761726
$ws.Range("I2").Value = -2938364.18659807
$ws.Range("J2").Value = 0.958038616411136
$ws.Range("K2").Value = 0.999569340352638
$ws.Range("L2").Value = -0.956629205321652
$ws.Range("M2").Value = 51.191570113583
$ws.Range("N2").Value = 2922.19331181278
$ws.Range("O2").Value = -2938364.18659807

# Row 3 (Group = 16)
$ws.Range("B3").Value = 0.958630856791998
$ws.Range("C3").Value = 208.691221224262
$ws.Range("D3").Value = -0.99186416443196
$ws.Range("E3").Value = -0.960485486419506
$ws.Range("F3").Value = 0.988893602938317
$ws.Range("G3").Value = -13.9503254317761
$ws.Range("H3").Value = -739.100551819794
$ws.Range("I3").Value = 3.93479185101298
$ws.Range("J3").Value = -0.959852434254271
$ws.Range("K3").Value = -0.858674536582941
$ws.Range("L3").Value = 0.988893602938317
$ws.Range("M3").Value = -0.057396432285014
$ws.Range("N3").Value = -2.80924484210361
$ws.Range("O3").Value = 3.93479185101298

# Row 5 (Group = 18)
$ws.Range("B5").Value = 0.959724684595245
$ws.Range("C5").Value = 210.937762444454
$ws.Range("D5").Value = -0.992003284461223
$ws.Range("E5").Value = -0.962931925400802
$ws.Range("F5").Value = 0.988964243022986
$ws.Range("G5").Value = -14.326424447208
$ws.Range("H5").Value = -760.853197313738
$ws.Range("I5").Value = 3.97290041025702
$ws.Range("J5").Value = -0.994758130049207
$ws.Range("K5").Value = -0.945219214524355
$ws.Range("L5").Value = 0.988964243022986
$ws.Range("M5").Value = -0.0872622655899818
$ws.Range("N5").Value = -4.53651247098088
$ws.Range("O5").Value = 3.97290041025702
